$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setup")
$ws.Rows("15:16").Delete()
$ws.Activate()
